$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) — row 3 -> F3: 3027 -> 3028, row 5 -> F5: 340 -> 370
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3028
$ws1.Range("F5").Value = 370

# Sheet "全部类型" (All Types) — same updates mirrored here
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F3").Value = 3028
$ws2.Range("F5").Value = 370
